# "Modifying variable name from suiteDxls to watchlistXls"
# The "A Suite" / IAM test-suite row becomes the "IAM" row and every
# Runmode in the Watchlist test suite is flipped to "Y".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Suite")

# Row 2 (TSID column) used to read "A Suite" - it now identifies the IAM suite.
$ws.Range("A2").Value = "IAM"

# Every row's Runmode now runs ("Y") instead of being skipped ("N").
foreach ($r in 2..7) {
    $ws.Cells.Item($r, 3).Value = "Y"
}

# Selection moves to the Runmode column that was just updated.
$ws.Range("C2:C7").Select()
